# Appointment reminder form: blood-draw related edits
#  - Drop the "blood_draw" note row from the survey sheet
#  - Convert the "upcoming" reminder note into a select_one question
#    that points at a new "soon" choice list
#  - Add the "soon" choice list (soon_noted / soon_keep) to the choices sheet

$wb  = $excel.ActiveWorkbook
$survey  = $wb.Worksheets.Item(1)
$choices = $wb.Worksheets.Item(2)

# 1. Remove the blood_draw note row (row 36) entirely; everything below
#    shifts up by one row automatically (formats included).
$survey.Rows.Item(36).Delete() | Out-Null

# 2. The old "upcoming" note (now row 38) becomes a select_one question,
#    and its label text is simplified (the two follow-up options move to
#    the new choice list instead of being embedded in the note text).
$survey.Range("A38").Value = "select_one soon"
$survey.Range("C38").Value = "Reminder! Client has an appointment soon.`n"

# 3. Add the "soon" choice list to the choices sheet.
$choices.Range("A2").Value = "soon"
$choices.Range("B2").Value = "soon_noted"
$choices.Range("C2").Value = "Noted, I will follow-up as needed. Delete this Task"
$choices.Range("A3").Value = "soon"
$choices.Range("B3").Value = "soon_keep"
$choices.Range("C3").Value = "Keep this reminder in my Task List"

# Match the formatting already used for this kind of row (copy it from the
# "follow-up" row that already carries the right style) onto the two new
# choice rows.
$survey.Range("A39").Copy() | Out-Null
$choices.Range("A2:C3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Widen the label column on the choices sheet to fit the new text.
$choices.Columns.Item(3).ColumnWidth = 36.3
